# Simulated Wild Card round and logged it
# Updates cumulative Rushing and Receiving stats with the results of the
# newly-simulated Wild Card round game, and adds a new Rushing row for
# D.Knox (his first logged carries of the season).

$wb = $excel.ActiveWorkbook

function Set-Row {
    param($ws, $row, $startCol, $values)
    $col = $startCol
    foreach ($val in $values) {
        $ws.Cells.Item($row, $col).Value = $val
        $col = $col + 1
    }
}

# ---------------------------------------------------------------------
# "Rushing" sheet (columns: A=idx, B=Name, C=1DATT, D=2DATT, E=3DATT, F=RZATT)
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

Set-Row $rushing 2 3 @(36, 38, 40, 33)     # J.Allen
Set-Row $rushing 3 3 @(109, 79, 7, 44)     # D.Singletary
Set-Row $rushing 4 3 @(50, 34, 8, 27)      # Z.Moss
Set-Row $rushing 5 3 @(12, 11, 3, 3)       # M.Breida (unchanged)
Set-Row $rushing 6 3 @(0, 1, 2, 2)         # R.Gilliam (unchanged)
Set-Row $rushing 7 3 @(2, 1, 0, 0)         # E.Sanders (unchanged)
Set-Row $rushing 8 3 @(5, 6, 2, 5)         # I.McKenzie

# New row for D.Knox, who logged his first carries this round
$rushing.Range("A9").Value = 7
$rushing.Range("B9").Value = "D.Knox"
Set-Row $rushing 9 3 @(1, 0, 0, 0)

# Match the formatting used by the rest of column A (bold/centered/bordered)
$rushing.Range("A8").Copy()
$rushing.Range("A9").PasteSpecial(-4122)
$rushing.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# "Receiving" sheet (columns: C=Short Target, D=Short Comp, E=Deep Target,
#  F=Deep Comp, G=RZ Target, H=RZ Comp)
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

Set-Row $receiving 2 3 @(51, 41, 1, 0, 5, 5)        # D.Singletary
Set-Row $receiving 3 3 @(30, 21, 1, 1, 6, 4)        # Z.Moss
Set-Row $receiving 4 3 @(8, 6, 1, 1, 3, 3)          # M.Breida (unchanged)
Set-Row $receiving 5 3 @(3, 2, 1, 0, 1, 0)          # R.Gilliam
Set-Row $receiving 6 3 @(123, 88, 34, 12, 28, 16)   # S.Diggs
Set-Row $receiving 7 3 @(47, 32, 28, 13, 8, 4)      # E.Sanders
Set-Row $receiving 8 3 @(104, 78, 10, 5, 15, 9)     # C.Beasley
Set-Row $receiving 9 3 @(39, 22, 25, 14, 18, 10)    # G.Davis
Set-Row $receiving 10 3 @(15, 12, 2, 0, 3, 3)       # I.McKenzie
Set-Row $receiving 11 3 @(2, 1, 1, 0, 1, 0)         # J.Kumerow
Set-Row $receiving 12 3 @(54, 43, 21, 14, 17, 12)   # D.Knox
Set-Row $receiving 13 3 @(10, 8, 0, 0, 1, 1)        # T.Sweeney (unchanged)
